$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of trade/prediction data appended below the existing 3 rows.
$rows = @(
    @{ Row = 4;  Date = 42628.834120370368; Score = -16; Buy = 1.66; PEG = 0.0969; T = 4.57; U = 4.5999999999999996; V = 2.2799999999999998 },
    @{ Row = 5;  Date = 42628.837546296294; Score = -17; Buy = 1.66; PEG = 0.0969; T = 4.57; U = 4.5999999999999996; V = 2.2799999999999998 },
    @{ Row = 6;  Date = 42628.838750000003; Score = 15;  Buy = 1.66; PEG = 0.0969; T = 4.57; U = 4.5999999999999996; V = 2.2799999999999998 },
    @{ Row = 7;  Date = 42628.840949074074; Score = -10; Buy = 1.66; PEG = 0.0969; T = 4.57; U = 4.5999999999999996; V = 2.2799999999999998 }
)

foreach ($r in $rows) {
    $i = $r.Row

    $ws.Range("A$i").Value = $r.Date
    $ws.Range("A$i").NumberFormat = "m/d/yy h:mm"

    $ws.Range("B$i").Value = $r.Score
    $ws.Range("C$i").Value = "buy"

    # Columns D through N: zeros
    for ($col = 4; $col -le 14; $col++) {
        $ws.Cells.Item($i, $col).Value = 0
    }

    $ws.Range("O$i").Value = 0
    $ws.Range("P$i").Value = "Random"
    $ws.Range("Q$i").Value = 0

    $ws.Range("R$i").Value = $r.Buy
    $ws.Range("S$i").Value = $r.PEG
    $ws.Range("S$i").NumberFormat = "0.00%"
    $ws.Range("T$i").Value = $r.T
    $ws.Range("U$i").Value = $r.U
    $ws.Range("V$i").Value = $r.V
    $ws.Range("W$i").Value = 0
}
